$d = $word.ActiveDocument

# Position right at the very end of the document (after the last paragraph's
# mark), so the new content is appended as a brand-new paragraph.
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/part" pkg:contentType="application/xml">' +
      '<pkg:xmlData>' +
        '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:pPr>' +
            '<w:pStyle w:val="a3"/>' +
            '<w:numPr>' +
              '<w:ilvl w:val="0"/>' +
              '<w:numId w:val="1"/>' +
            '</w:numPr>' +
          '</w:pPr>' +
          '<w:r>' +
            '<w:t>Добавить модуль защиты аккумулятора от переразряда</w:t>' +
          '</w:r>' +
        '</w:p>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertionPoint.InsertXML($newParaXml)
